$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.670.22"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.848.00"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4265"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3629"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07308"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8746"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "1.905.31"
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.318"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.508"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06906"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009011"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "27.686.59"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.948"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value = "2.143.90"
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.971"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("E29").Value = "  +9.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.261"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.858"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08925"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7597"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.968"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.521"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.098"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05380"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.093"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01931"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.822"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5068"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1653"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.761"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.330"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06551"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4667"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.616"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.12%  "
